# Weekly update: insert two new daily records at the top of the "Papa" data
# block (rows 879-880), shifting all subsequent rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 879; this pushes the
# existing rows 879.. down to 881.. (matching dimension growth R951 -> R953)
$ws.Rows("879:880").Insert()

# New row 879
$ws.Cells.Item(879, 1).Value = 5
$ws.Cells.Item(879, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(879, 3).Value = "Maule"
$ws.Cells.Item(879, 4).Value = 45223
$ws.Cells.Item(879, 5).Value = 7
$ws.Cells.Item(879, 6).Value = 100114001
$ws.Cells.Item(879, 7).Value = "Papa"
$ws.Cells.Item(879, 8).Value = "Rodeo"
$ws.Cells.Item(879, 9).Value = "1a (guarda lavada)"
$ws.Cells.Item(879, 10).Value = 1500
$ws.Cells.Item(879, 11).Value = 28000
$ws.Cells.Item(879, 12).Value = 28000
$ws.Cells.Item(879, 13).Value = 28000
$ws.Cells.Item(879, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(879, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(879, 16).Value = 1120
$ws.Cells.Item(879, 17).Value = 25
$ws.Cells.Item(879, 18).Value = "Hortaliza"

# New row 880
$ws.Cells.Item(880, 1).Value = 5
$ws.Cells.Item(880, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(880, 3).Value = "Maule"
$ws.Cells.Item(880, 4).Value = 45223
$ws.Cells.Item(880, 5).Value = 7
$ws.Cells.Item(880, 6).Value = 100114001
$ws.Cells.Item(880, 7).Value = "Papa"
$ws.Cells.Item(880, 8).Value = "Rodeo"
$ws.Cells.Item(880, 9).Value = "1a (guarda)"
$ws.Cells.Item(880, 10).Value = 800
$ws.Cells.Item(880, 11).Value = 26000
$ws.Cells.Item(880, 12).Value = 26000
$ws.Cells.Item(880, 13).Value = 26000
$ws.Cells.Item(880, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(880, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(880, 16).Value = 1040
$ws.Cells.Item(880, 17).Value = 25
$ws.Cells.Item(880, 18).Value = "Hortaliza"
